$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) updates for the crypto list refresh.
# D-column values are numeric-looking text (e.g. "26.384.36"); they must
# stay stored as text, matching the original inlineStr cells, so we
# temporarily force a text number format before assignment and then
# restore the default "Normal" style so no stray formatting is left behind.

$d = $ws.Range("D2")
$d.NumberFormat = "@"
$d.Value = "26.384.36"
$d.Style = "Normal"
$ws.Range("E2").Value = "  -0.44%  "

$d = $ws.Range("D3")
$d.NumberFormat = "@"
$d.Value = "1.717.67"
$d.Style = "Normal"
$ws.Range("E3").Value = "  -0.95%  "

$d = $ws.Range("D4")
$d.NumberFormat = "@"
$d.Value = "0.9977"
$d.Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "

$d = $ws.Range("D5")
$d.NumberFormat = "@"
$d.Value = "241.72"
$d.Style = "Normal"
$ws.Range("E5").Value = "  -2.08%  "

$d = $ws.Range("D6")
$d.NumberFormat = "@"
$d.Value = "0.9983"
$d.Style = "Normal"
$ws.Range("E6").Value = "  -0.19%  "

$d = $ws.Range("D7")
$d.NumberFormat = "@"
$d.Value = "0.4864"
$d.Style = "Normal"
$ws.Range("E7").Value = "  -0.47%  "

$d = $ws.Range("D8")
$d.NumberFormat = "@"
$d.Value = "0.2585"
$d.Style = "Normal"
$ws.Range("E8").Value = "  -3.05%  "

$d = $ws.Range("D9")
$d.NumberFormat = "@"
$d.Value = "0.06167"
$d.Style = "Normal"
$ws.Range("E9").Value = "  -2.21%  "

$d = $ws.Range("D10")
$d.NumberFormat = "@"
$d.Value = "1.725.85"
$d.Style = "Normal"
$ws.Range("E10").Value = "  -0.42%  "

$d = $ws.Range("D11")
$d.NumberFormat = "@"
$d.Value = "0.06949"
$d.Style = "Normal"
$ws.Range("E11").Value = "  -1.09%  "

$d = $ws.Range("D12")
$d.NumberFormat = "@"
$d.Value = "15.49"
$d.Style = "Normal"
$ws.Range("E12").Value = "  -1.21%  "

$d = $ws.Range("D13")
$d.NumberFormat = "@"
$d.Value = "4.490"
$d.Style = "Normal"
$ws.Range("E13").Value = "  -2.26%  "

$d = $ws.Range("D14")
$d.NumberFormat = "@"
$d.Value = "0.5974"
$d.Style = "Normal"
$ws.Range("E14").Value = "  -1.80%  "

$d = $ws.Range("D15")
$d.NumberFormat = "@"
$d.Value = "76.60"
$d.Style = "Normal"
$ws.Range("E15").Value = "  -0.96%  "

$d = $ws.Range("D16")
$d.NumberFormat = "@"
$d.Value = "0.9978"
$d.Style = "Normal"
$ws.Range("E16").Value = "  -0.23%  "

$d = $ws.Range("D17")
$d.NumberFormat = "@"
$d.Value = "26.379.66"
$d.Style = "Normal"
$ws.Range("E17").Value = "  -0.41%  "

$d = $ws.Range("D18")
$d.NumberFormat = "@"
$d.Value = "0.9983"
$d.Style = "Normal"
$ws.Range("E18").Value = "  -0.18%  "

$d = $ws.Range("D19")
$d.NumberFormat = "@"
$d.Value = "0.000007102"
$d.Style = "Normal"
$ws.Range("E19").Value = "  -5.30%  "

$d = $ws.Range("D20")
$d.NumberFormat = "@"
$d.Value = "11.24"
$d.Style = "Normal"
$ws.Range("E20").Value = "  -2.38%  "

$d = $ws.Range("D21")
$d.NumberFormat = "@"
$d.Value = "1.944.89"
$d.Style = "Normal"
$ws.Range("E21").Value = "  -0.38%  "

$d = $ws.Range("D22")
$d.NumberFormat = "@"
$d.Value = "4.415"
$d.Style = "Normal"
$ws.Range("E22").Value = "  -3.53%  "

$d = $ws.Range("D23")
$d.NumberFormat = "@"
$d.Value = "8.451"
$d.Style = "Normal"
$ws.Range("E23").Value = "  -2.90%  "

$d = $ws.Range("D24")
$d.NumberFormat = "@"
$d.Value = "5.070"
$d.Style = "Normal"
$ws.Range("E24").Value = "  -2.98%  "

$d = $ws.Range("D25")
$d.NumberFormat = "@"
$d.Value = "136.89"
$d.Style = "Normal"
$ws.Range("E25").Value = "  -2.70%  "

$d = $ws.Range("D26")
$d.NumberFormat = "@"
$d.Value = "15.22"
$d.Style = "Normal"
$ws.Range("E26").Value = "  -1.41%  "

$d = $ws.Range("D27")
$d.NumberFormat = "@"
$d.Value = "1.401"
$d.Style = "Normal"
$ws.Range("E27").Value = "  -1.06%  "

$d = $ws.Range("D28")
$d.NumberFormat = "@"
$d.Value = "1.735"
$d.Style = "Normal"
$ws.Range("E28").Value = "  -1.80%  "

$d = $ws.Range("D29")
$d.NumberFormat = "@"
$d.Value = "105.96"
$d.Style = "Normal"
$ws.Range("E29").Value = "  -1.84%  "

$d = $ws.Range("D30")
$d.NumberFormat = "@"
$d.Value = "3.881"
$d.Style = "Normal"
$ws.Range("E30").Value = "  -3.75%  "

$d = $ws.Range("D31")
$d.NumberFormat = "@"
$d.Value = "0.07968"
$d.Style = "Normal"
$ws.Range("E31").Value = "  -0.47%  "

$d = $ws.Range("D32")
$d.NumberFormat = "@"
$d.Value = "3.611"
$d.Style = "Normal"
$ws.Range("E32").Value = "  -2.74%  "

$d = $ws.Range("D33")
$d.NumberFormat = "@"
$d.Value = "0.04445"
$d.Style = "Normal"
$ws.Range("E33").Value = "  -2.83%  "

$d = $ws.Range("D34")
$d.NumberFormat = "@"
$d.Value = "2.601"
$d.Style = "Normal"
$ws.Range("E34").Value = "  -0.43%  "

$d = $ws.Range("D35")
$d.NumberFormat = "@"
$d.Value = "0.9952"
$d.Style = "Normal"
$ws.Range("E35").Value = "  -1.23%  "

$d = $ws.Range("D36")
$d.NumberFormat = "@"
$d.Value = "0.6190"
$d.Style = "Normal"
$ws.Range("E36").Value = "  -2.53%  "

$d = $ws.Range("D37")
$d.NumberFormat = "@"
$d.Value = "0.9343"
$d.Style = "Normal"
$ws.Range("E37").Value = "  +4.54%  "

$d = $ws.Range("D38")
$d.NumberFormat = "@"
$d.Value = "1.970"
$d.Style = "Normal"
$ws.Range("E38").Value = "  -2.10%  "

$d = $ws.Range("D39")
$d.NumberFormat = "@"
$d.Value = "2.375"
$d.Style = "Normal"
$ws.Range("E39").Value = "  -0.92%  "

$ws.Range("E40").Value = "  -0.47%  "

$d = $ws.Range("D41")
$d.NumberFormat = "@"
$d.Value = "0.01476"
$d.Style = "Normal"
$ws.Range("E41").Value = "  -1.80%  "

$d = $ws.Range("D42")
$d.NumberFormat = "@"
$d.Value = "99.16"
$d.Style = "Normal"
$ws.Range("E42").Value = "  -2.53%  "

$d = $ws.Range("D43")
$d.NumberFormat = "@"
$d.Value = "5.455"
$d.Style = "Normal"
$ws.Range("E43").Value = "  +0.95%  "

$d = $ws.Range("D44")
$d.NumberFormat = "@"
$d.Value = "0.3816"
$d.Style = "Normal"
$ws.Range("E44").Value = "  -1.70%  "

$d = $ws.Range("D45")
$d.NumberFormat = "@"
$d.Value = "6.829"
$d.Style = "Normal"
$ws.Range("E45").Value = "  -1.22%  "

$d = $ws.Range("D46")
$d.NumberFormat = "@"
$d.Value = "0.1151"
$d.Style = "Normal"
$ws.Range("E46").Value = "  -2.64%  "

$d = $ws.Range("D47")
$d.NumberFormat = "@"
$d.Value = "0.05353"
$d.Style = "Normal"
$ws.Range("E47").Value = "  -0.70%  "

$d = $ws.Range("D48")
$d.NumberFormat = "@"
$d.Value = "30.40"
$d.Style = "Normal"
$ws.Range("E48").Value = "  -0.21%  "

$d = $ws.Range("D49")
$d.NumberFormat = "@"
$d.Value = "7.709"
$d.Style = "Normal"
$ws.Range("E49").Value = "  -0.50%  "

$d = $ws.Range("D50")
$d.NumberFormat = "@"
$d.Value = "51.22"
$d.Style = "Normal"
$ws.Range("E50").Value = "  -0.94%  "

$d = $ws.Range("D51")
$d.NumberFormat = "@"
$d.Value = "1.218"
$d.Style = "Normal"
$ws.Range("E51").Value = "  -3.16%  "

